# Rename: "Demo App zur DWX"
#
# Changes applied to the second paragraph:
#   "Die Anwendung zeigt den die Verwendung von Apache Cordova ..."
#   -> "Die Anwendung zeigt die Verwendung von Apache Cordova ..."
# (drops the superfluous "den", and the grammar-check proofErr markers
# that used to wrap "Apache" disappear with it)
#
# The "_GoBack" bookmark, which used to wrap the entire second
# paragraph, is collapsed to an empty bookmark placed right after
# "Die Anwendung zeigt " (where Word leaves it after the last edit).

$d = $word.ActiveDocument

# Locate "den die Verwendung von Apache " (including the trailing
# space) using Find so we don't have to hard-code character offsets.
$target = $d.Content
$target.Find.Execute("den die Verwendung von Apache ") | Out-Null

# Remember where the replacement text starts; this is also where the
# "_GoBack" bookmark needs to move to.
$insertStart = $target.Start

# Remove the matched text (this also removes the now-enclosed
# proofErr gramStart/gramEnd markers around "Apache").
$target.Delete()

# Insert the corrected wording in its place.
$insertRange = $d.Range($insertStart, $insertStart)
$insertRange.InsertAfter("die Verwendung von Apache ")

# Move the "_GoBack" bookmark to an empty range right before the
# newly inserted text; re-adding a bookmark with the same name moves
# it (removing the old bookmarkStart/bookmarkEnd pair).
$bookmarkRange = $d.Range($insertStart, $insertStart)
$d.Bookmarks.Add("_GoBack", $bookmarkRange) | Out-Null
